# Auto-generated edit script: apply updated Market Board values to the Leve profit sheets.
# Source: scheduled runner diff (see commit message) -- values only, no formula/style changes.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 788.5
$ws.Range("I18").Value = 788.5
$ws.Range("K18").Value = 788.5
$ws.Range("M18").Value = -504.5
$ws.Range("H21").Value = 9902.666999999999
$ws.Range("I21").Value = 9902.666999999999
$ws.Range("K21").Value = 9902.666999999999
$ws.Range("M21").Value = -9434.666999999999
$ws.Range("H23").Value = 9902.666999999999
$ws.Range("I23").Value = 9902.666999999999
$ws.Range("K23").Value = 9902.666999999999
$ws.Range("M23").Value = -9668.666999999999
$ws.Range("H40").Value = 4804.619
$ws.Range("I40").Value = 3825.8462
$ws.Range("K40").Value = 3825.8462
$ws.Range("M40").Value = -3650.8462
$ws.Range("H43").Value = 2173.1428
$ws.Range("I43").Value = 1500.5
$ws.Range("K43").Value = 1500.5
$ws.Range("M43").Value = -1431.5
$ws.Range("H58").Value = 978
$ws.Range("I58").Value = 71.75
$ws.Range("J58").Value = 1582.1666
$ws.Range("K58").Value = 215.25
$ws.Range("L58").Value = 4746.4998
$ws.Range("M58").Value = -65.25
$ws.Range("N58").Value = -5046.4998
$ws.Range("H64").Value = 9248
$ws.Range("I64").Value = 8500.5
$ws.Range("K64").Value = 8500.5
$ws.Range("M64").Value = -8252.5
$ws.Range("H67").Value = 9248
$ws.Range("I67").Value = 8500.5
$ws.Range("K67").Value = 8500.5
$ws.Range("M67").Value = -7642.5
$ws.Range("H92").Value = 3983.1667
$ws.Range("I92").Value = 2633.3333
$ws.Range("J92").Value = 5333
$ws.Range("K92").Value = 2633.3333
$ws.Range("L92").Value = 5333
$ws.Range("M92").Value = -1385.3333
$ws.Range("N92").Value = -7829
$ws.Range("H137").Value = 2870.9697
$ws.Range("J137").Value = 4207.2144
$ws.Range("L137").Value = 12621.6432
$ws.Range("N137").Value = -17721.6432
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2217.88
$ws.Range("I2").Value = 1550.8948
$ws.Range("J2").Value = 4330
$ws.Range("K2").Value = 1550.8948
$ws.Range("L2").Value = 4330
$ws.Range("M2").Value = -1437.8948
$ws.Range("N2").Value = -4556
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("M29").ClearContents()
$ws.Range("H32").Value = 4216.769
$ws.Range("I32").Value = 4216.769
$ws.Range("K32").Value = 4216.769
$ws.Range("M32").Value = -3929.769
$ws.Range("H45").Value = 2828.4583
$ws.Range("I45").Value = 2375.3809
$ws.Range("K45").Value = 2375.3809
$ws.Range("M45").Value = -1998.3809
$ws.Range("H97").Value = 567.8333
$ws.Range("I97").Value = 542.41174
$ws.Range("J97").Value = 1000
$ws.Range("K97").Value = 542.41174
$ws.Range("L97").Value = 1000
$ws.Range("M97").Value = -46.41174000000001
$ws.Range("N97").Value = -1992
$ws.Range("H110").Value = 3135.6667
$ws.Range("J110").Value = 5599.4
$ws.Range("L110").Value = 5599.4
$ws.Range("N110").Value = -9689.4
$ws.Range("H116").Value = 2217.88
$ws.Range("I116").Value = 1550.8948
$ws.Range("J116").Value = 4330
$ws.Range("K116").Value = 1550.8948
$ws.Range("L116").Value = 4330
$ws.Range("M116").Value = 743.1052
$ws.Range("N116").Value = -8918
$ws.Range("H122").Value = 1532.4166
$ws.Range("I122").Value = 1532.4166
$ws.Range("K122").Value = 4597.2498
$ws.Range("M122").Value = -2147.2498
$ws.Range("H129").Value = 0
$ws.Range("J129").Value = 0
$ws.Range("L129").Value = 0
$ws.Range("N129").ClearContents()
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2217.88
$ws.Range("I3").Value = 1550.8948
$ws.Range("J3").Value = 4330
$ws.Range("K3").Value = 1550.8948
$ws.Range("L3").Value = 4330
$ws.Range("M3").Value = -1436.8948
$ws.Range("N3").Value = -4558
$ws.Range("H86").Value = 5274.2104
$ws.Range("I86").Value = 4301.0713
$ws.Range("K86").Value = 4301.0713
$ws.Range("M86").Value = -3178.0713
$ws.Range("H89").Value = 5274.2104
$ws.Range("I89").Value = 4301.0713
$ws.Range("K89").Value = 21505.3565
$ws.Range("M89").Value = -15889.3565
$ws.Range("H94").Value = 612.5
$ws.Range("I94").Value = 557.7143
$ws.Range("K94").Value = 557.7143
$ws.Range("M94").Value = -106.7143
$ws.Range("H134").Value = 2875.4722
$ws.Range("I134").Value = 2671.9143
$ws.Range("K134").Value = 8015.742899999999
$ws.Range("M134").Value = -5480.742899999999
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1768.25
$ws.Range("I16").Value = 1768.25
$ws.Range("K16").Value = 1768.25
$ws.Range("M16").Value = -1481.25
$ws.Range("H22").Value = 2656
$ws.Range("I22").Value = 524.5
$ws.Range("J22").Value = 3366.5
$ws.Range("K22").Value = 524.5
$ws.Range("L22").Value = 3366.5
$ws.Range("M22").Value = -174.5
$ws.Range("N22").Value = -4066.5
$ws.Range("H31").Value = 5678.65
$ws.Range("I31").Value = 2345.9546
$ws.Range("K31").Value = 2345.9546
$ws.Range("M31").Value = -2050.9546
$ws.Range("H34").Value = 5678.65
$ws.Range("I34").Value = 2345.9546
$ws.Range("K34").Value = 2345.9546
$ws.Range("M34").Value = -2143.9546
$ws.Range("H58").Value = 4920.654
$ws.Range("I58").Value = 2593.6316
$ws.Range("K58").Value = 2593.6316
$ws.Range("M58").Value = -2390.6316
$ws.Range("H113").Value = 1768.25
$ws.Range("I113").Value = 1768.25
$ws.Range("K113").Value = 1768.25
$ws.Range("M113").Value = 401.75
$ws.Range("H136").Value = 4920.654
$ws.Range("I136").Value = 2593.6316
$ws.Range("K136").Value = 7780.8948
$ws.Range("M136").Value = -5230.8948
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H114").Value = 3477.5833
$ws.Range("J114").Value = 3703.4443
$ws.Range("L114").Value = 11110.3329
$ws.Range("N114").Value = -17618.3329
$ws.Range("H119").Value = 0
$ws.Range("I119").Value = 0
$ws.Range("K119").Value = 0
$ws.Range("M119").ClearContents()
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 10750.25
$ws.Range("J21").Value = 14500
$ws.Range("L21").Value = 14500
$ws.Range("N21").Value = -14846
$ws.Range("H30").Value = 10750.25
$ws.Range("J30").Value = 14500
$ws.Range("L30").Value = 14500
$ws.Range("N30").Value = -14710
$ws.Range("H47").Value = 29899
$ws.Range("J47").Value = 29899
$ws.Range("L47").Value = 29899
$ws.Range("N47").Value = -31035
$ws.Range("H55").Value = 6789.857
$ws.Range("J55").Value = 8349.5
$ws.Range("L55").Value = 8349.5
$ws.Range("N55").Value = -9003.5
$ws.Range("H113").Value = 9790.833000000001
$ws.Range("I113").Value = 9375
$ws.Range("K113").Value = 9375
$ws.Range("M113").Value = -7205
$ws.Range("H126").Value = 3409.2942
$ws.Range("I126").Value = 3130.5334
$ws.Range("K126").Value = 9391.600199999999
$ws.Range("M126").Value = -6921.600199999999
$ws.Range("H132").Value = 57223.1
$ws.Range("I132").Value = 62247.945
$ws.Range("K132").Value = 186743.835
$ws.Range("M132").Value = -184213.835
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H12").Value = 524.5
$ws.Range("J12").Value = 524.5
$ws.Range("L12").Value = 524.5
$ws.Range("N12").Value = -864.5
$ws.Range("H26").Value = 63010
$ws.Range("J26").Value = 63010
$ws.Range("L26").Value = 63010
$ws.Range("N26").Value = -63600
$ws.Range("H46").Value = 9977.272000000001
$ws.Range("J46").Value = 8222.223
$ws.Range("L46").Value = 8222.223
$ws.Range("N46").Value = -8598.223
$ws.Range("H61").Value = 2964.48
$ws.Range("I61").Value = 2277.818
$ws.Range("J61").Value = 8000
$ws.Range("K61").Value = 2277.818
$ws.Range("L61").Value = 8000
$ws.Range("M61").Value = -2075.818
$ws.Range("N61").Value = -8404
$ws.Range("H113").Value = 2964.48
$ws.Range("I113").Value = 2277.818
$ws.Range("J113").Value = 8000
$ws.Range("K113").Value = 2277.818
$ws.Range("L113").Value = 8000
$ws.Range("M113").Value = -107.8180000000002
$ws.Range("N113").Value = -12340
$ws.Range("H132").Value = 5995
$ws.Range("I132").Value = 5995
$ws.Range("K132").Value = 17985
$ws.Range("M132").Value = -15455
$ws.Range("H136").Value = 5282.5713
$ws.Range("I136").Value = 4995.8
$ws.Range("K136").Value = 14987.4
$ws.Range("M136").Value = -12437.4
